# Fixed a bad reference in data export in aL extract.
# Added line in ellipsoid to fix bad segmentation error in reflexed plane systems.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New MgZn2 data rows (rows 7-15), mirroring the existing NaCl/CsCl block layout.
$data = @(
    @("MgZn2", 0.92,  8,   -4,  3),
    @("MgZn2", 0.852, 14,  -10, -1),
    @("MgZn2", 0.815, 11,  -8,  -3),
    @("MgZn2", 0.782, 0,   -1,  -4),
    @("MgZn2", 0.69,  -4,  20,  10),
    @("MgZn2", 0.758, -19, 5,   -7),
    @("MgZn2", 0.692, -24, 30,  9),
    @("MgZn2", 0.797, 7,   -19, -9),
    @("MgZn2", 0.673, -16, 34,  13)
)

$row = 7
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# Widen column B slightly so the new fractional values read cleanly.
$ws.Columns("B").ColumnWidth = 12.43

# Reselect/scroll to frame the newly appended block and zoom out a notch.
$null = $ws.Range("A7:E15").Select()
$excel.ActiveWindow.Zoom = 90
